$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: task text was a typo ("Iplement classes...") - replace with corrected task
# and update the estimated hours for it (3 -> 6)
$ws.Range("A2").Value = "Implement calsses according to UML-class diagram without virus"
$ws.Range("B2").Value = 6

# Remaining tasks keep their text but get new "hours per person" estimates
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 15
$ws.Range("B6").Value = 15
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 18

# Move the active selection
[void]$ws.Range("A13").Select()
